# Circle Language Spec Plan: Set font to Calibri for non-heading text.
#
# This reproduces (as closely as the Word object model allows):
#   1. The stray "_GoBack" bookmark is moved from right after "2008-03 "
#      (before "Symbol = Creator" in the title) down into the body, landing
#      mid-word inside "Finish the article " -> "Finis" | "h the article ".
#      (This is simply Word re-stamping its "last edit position" bookmark
#      after the real edit below was made - net effect: delete the old
#      bookmark, add it back at the new cursor position.)
#   2. The "Normal" style's font is changed from Tahoma 12pt to Calibri 11pt
#      (the actual content edit the commit message is about).

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark -------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $d.Content
$null = $target.Find.Execute("Finis", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$splitPoint = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# --- 2. Re-style "Normal" to Calibri 11pt (non-heading body text font) --

$normal = $d.Styles("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 11
